# Insert a new row at position 206, shifting existing rows 206..255 down to 207..256,
# then populate the newly inserted row 206 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 206.
$ws.Rows.Item(206).Insert()

# Populate the new row 206 with the new data record.
$ws.Range("A206").Value2 = 6
$ws.Range("B206").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C206").Value2 = "Metropolitana"
$ws.Range("D206").Value2 = 44511
$ws.Range("E206").Value2 = 13
$ws.Range("F206").Value2 = 100112032
$ws.Range("G206").Value2 = "Zapallo italiano"
$ws.Range("H206").Value2 = "Sin especificar"
$ws.Range("I206").Value2 = "Primera"
$ws.Range("J206").Value2 = 400
$ws.Range("K206").Value2 = 5000
$ws.Range("L206").Value2 = 6000
$ws.Range("M206").Value2 = 5575
$ws.Range("N206").Value2 = "`$/caja 50 unidades"
$ws.Range("O206").Value2 = "Región de O'Higgins"
$ws.Range("P206").Value2 = 112
$ws.Range("Q206").Value2 = 50
$ws.Range("R206").Value2 = "Hortaliza"
